# Applies the Yarpiz PSO table edit:
#  - rename header "Gen" -> "MaxFES"
#  - delete the "Run 50" column (last run column), shifting the
#    trailing "Mean" column one place to the left
#  - replace the MaxFES (column A) values for rows 2-14 with the new
#    fractional progress values
#  - replace the recalculated Mean column values for rows 2-14

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "Run 50" column (AZ); cells to its right (the old
# "Mean" column, BA) shift left to become the new AZ column.
$ws.Range("AZ1:AZ14").Delete()

# Header rename.
$ws.Range("A1").Value = "MaxFES"

# Column A (MaxFES) data values.
$ws.Range("A2").Value = 0
$ws.Range("A3").Value = 0.001
$ws.Range("A4").Value = 0.01
$ws.Range("A5").Value = 0.1
$ws.Range("A6").Value = 0.2
$ws.Range("A7").Value = 0.3
$ws.Range("A8").Value = 0.4
$ws.Range("A9").Value = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# Recalculated Mean column (now AZ after the shift) data values.
$ws.Range("AZ2").Value = 62.85581152
$ws.Range("AZ3").Value = 42.22975964
$ws.Range("AZ4").Value = 1.13330657
$ws.Range("AZ5").Value = 0.65165418
$ws.Range("AZ6").Value = 0.65165418
$ws.Range("AZ7").Value = 0.65165418
$ws.Range("AZ8").Value = 0.65165418
$ws.Range("AZ9").Value = 0.65165418
$ws.Range("AZ10").Value = 0.65165418
$ws.Range("AZ11").Value = 0.65165418
$ws.Range("AZ12").Value = 0.65165418
$ws.Range("AZ13").Value = 0.65165418
$ws.Range("AZ14").Value = 0.65165418
